$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @(
    'Última fecha'
    '2024-01-16 00:00:00'
    '2023-12-29 00:00:00'
    '2023-12-22 00:00:00'
    '2024-01-16 00:00:00'
    '2024-01-16 00:00:00'
    '2023-12-08 00:00:00'
    '2023-12-21 00:00:00'
    '2024-01-02 00:00:00'
    '2024-01-02 00:00:00'
    '2024-01-23 00:00:00'
    '2023-12-29 00:00:00'
    '2023-12-13 00:00:00'
    '2023-10-20 00:00:00'
    '2023-11-27 00:00:00'
    '2023-11-27 00:00:00'
    '2024-01-16 00:00:00'
    '2024-01-23 00:00:00'
    '2024-01-02 00:00:00'
    '2023-12-29 00:00:00'
    '2024-01-16 00:00:00'
    '2024-01-23 00:00:00'
    '2024-01-23 00:00:00'
    '2024-01-16 00:00:00'
    '2024-01-23 00:00:00'
    '2024-01-16 00:00:00'
    '2024-01-23 00:00:00'
    '2024-01-16 00:00:00'
    '2024-01-16 00:00:00'
    '2024-01-23 00:00:00'
    '2024-01-16 00:00:00'
    '2024-01-16 00:00:00'
    '2024-01-16 00:00:00'
    '2023-12-29 00:00:00'
    '2024-01-16 00:00:00'
    '2023-12-29 00:00:00'
    '2024-01-16 00:00:00'
    '2024-01-16 00:00:00'
    '2023-12-21 00:00:00'
    '2023-12-21 00:00:00'
    '2024-01-16 00:00:00'
    '2024-01-23 00:00:00'
    '2023-12-29 00:00:00'
    '2024-01-12 00:00:00'
    '2024-01-16 00:00:00'
    '2024-01-05 00:00:00'
    '2023-12-29 00:00:00'
    '2024-01-16 00:00:00'
    '2023-12-21 00:00:00'
    '2023-12-21 00:00:00'
    '2024-01-16 00:00:00'
    '2023-12-21 00:00:00'
    '2023-12-21 00:00:00'
    '2023-12-21 00:00:00'
    '2023-12-21 00:00:00'
    '2024-01-23 00:00:00'
    '2023-12-21 00:00:00'
    '2023-12-21 00:00:00'
    '2024-01-02 00:00:00'
    '2023-12-21 00:00:00'
    '2023-12-21 00:00:00'
    '2023-12-21 00:00:00'
    '2023-12-21 00:00:00'
    '2023-12-13 00:00:00'
    '2023-12-21 00:00:00'
    '2023-12-21 00:00:00'
    '2024-01-16 00:00:00'
    '2024-01-16 00:00:00'
    '2024-01-16 00:00:00'
    '2024-01-16 00:00:00'
    '2024-01-16 00:00:00'
    '2024-01-02 00:00:00'
    '2024-01-02 00:00:00'
    '2024-01-16 00:00:00'
    '2024-01-16 00:00:00'
    '2024-01-23 00:00:00'
    '2023-11-27 00:00:00'
    '2023-11-27 00:00:00'
    '2024-01-16 00:00:00'
    '2023-11-27 00:00:00'
    '2023-11-27 00:00:00'
    '2024-01-23 00:00:00'
    '2024-01-16 00:00:00'
    '2024-01-16 00:00:00'
    '2024-01-16 00:00:00'
    '2023-12-21 00:00:00'
    '2023-12-29 00:00:00'
    '2023-12-20 00:00:00'
    '2023-12-20 00:00:00'
    '2024-01-11 00:00:00'
    '2024-01-29 00:00:00'
    '2024-01-11 00:00:00'
)

for ($i = 0; $i -lt $gValues.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 7).Value = $gValues[$i]
}

